$d = $word.ActiveDocument

# --- Change 2 (done first): remove the old "_GoBack" bookmark sitting between
# "autore" / "s" and let those two runs merge back into a single "autores" run.
# (Word only allows one bookmark per name, so the old one has to go before we
# mint a new "_GoBack" bookmark elsewhere in the document.)
$oldBookmark = $d.Bookmarks("_GoBack")
$searchStart = $oldBookmark.Start - 10
$oldBookmark.Delete()

# Scope the Find to the tail of the document (well past the earlier, unrelated
# "...registros de autores." sentence) so only the split "autore"+"s" occurrence
# gets touched/merged — an unscoped Find("autores") would also "replace" the
# already-whole "autores" elsewhere and strip its rsid attribute as a side effect.
$scoped = $d.Range($searchStart, $d.Content.End)
$scoped.Find.Execute("autores", $true, $false, $false, $false, $false, $true, 1, $false, "autores", 2)

# --- Change 1: insert a fresh "_GoBack" bookmark at the very start of the document ---
# Bookmarks.Add collapses correctly for any non-zero collapsed Range, but a literal
# (0,0) Range is special-cased by the host as "whole document", so we insert a
# throw-away placeholder character first to get a safe non-zero anchor to bookmark,
# then remove the placeholder again.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("~")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range(0, 1).Delete()
